{"js": "// Remove the \"IT Support Intern\" text left in the \"Department Assigned\"\n// table cell, and remove the stray \"_GoBack\" bookmark left in the final\n// empty paragraph of the document body (both leftovers from the author's\n// last manual edit in Word).\n\nconst body = context.document.body;\n\n// --- 1. Remove the \"IT Support Intern\" run from the table cell --------\n// Search for the exact run text and blank it out in place; this removes\n// the <w:r> but keeps the owning (now-empty) paragraph and its <w:pPr>\n// formatting untouched, matching Word's own \"select text, press Delete\"\n// behavior.\nconst searchResults = body.search(\"IT Support Intern\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (const found of searchResults.items) {\n  found.insertText(\"\", \"Replace\");\n}\nawait context.sync();\n\n// --- 2. Remove the \"_GoBack\" bookmark ----------------------------------\nconst goBackRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBackRange.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBackRange.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Remove the \"IT Support Intern\" text left in the \"Department Assigned\"\n# table cell, and remove the stray \"_GoBack\" bookmark left in the final\n# empty paragraph of the document body (both leftovers from the author's\n# last manual edit in Word).\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the \"IT Support Intern\" run from the table cell --------\n# Find the exact run text and blank the matched range in place; this\n# removes the <w:r> but leaves the owning (now-empty) paragraph and its\n# <w:pPr> formatting untouched, matching Word's own \"select text, press\n# Delete\" behavior.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"IT Support Intern\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \"\"\n}\n\n# --- 2. Remove the \"_GoBack\" bookmark ----------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
